# feat: add 2022-Q1 data
#
# The existing "总计" (Total) sheet is renamed to "2022-Q1" and repurposed to
# hold the fund-level breakdown for the new quarter (matching the pattern used
# by the other per-quarter sheets). A brand-new "总计" sheet is appended after
# it, holding the same aggregate table as before plus a new row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Repurpose the old "总计" sheet -> "2022-Q1" fund-detail sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Grab a header row + index cell that already carry the bold/bordered style
# used across every other quarterly sheet, so the new sheet matches formatting.
$styleSrc = $wb.Worksheets.Item("2021-Q4")

# Wipe the old aggregate-table content/formatting before laying out the new
# fund-detail table.
$q1.Cells.Clear()

# Header row (B1:H1), copied (with formatting) from an existing sheet, then
# overwritten with this sheet's header text.
$styleSrc.Range("B1:H1").Copy($q1.Range("B1:H1"))
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Index column (A2:A3), copied (with formatting) for the styled numeric index.
$styleSrc.Range("A2").Copy($q1.Range("A2"))
$styleSrc.Range("A2").Copy($q1.Range("A3"))
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1

# Row 2 - 004351 / 汇丰晋信珠三角区域发展混合
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "004351"
$q1.Range("C2").Value = "汇丰晋信珠三角区域发展混合"
$q1.Range("D2").Value = "0.51"
$q1.Range("E2").Value = "93.92"
$q1.Range("F2").Value = "5.29"
$q1.Range("G2").Value = "0.0270"
$q1.Range("H2").Value = 3

# Row 3 - 002152 / 华宝核心优势灵活配置混合
$q1.Range("B3:G3").NumberFormat = "@"
$q1.Range("B3").Value = "002152"
$q1.Range("C3").Value = "华宝核心优势灵活配置混合"
$q1.Range("D3").Value = "0.45"
$q1.Range("E3").Value = "90.91"
$q1.Range("F3").Value = "3.48"
$q1.Range("G3").Value = "0.0157"
$q1.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 2) Add a fresh "总计" sheet after "2022-Q1" with the aggregate table
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Header row (B1:D1), copied (with formatting) for consistency.
$styleSrc.Range("B1:D1").Copy($total.Range("B1:D1"))
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# Styled index column A2:A7
$styleSrc.Range("A2").Copy($total.Range("A2:A7"))
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# Row 2 - new 2022-Q1 summary
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.04

# Row 3 - 2021-Q4 (previously row 2)
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.51

# Row 4 - 2021-Q3 (previously row 3)
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.03

# Row 5 - 2021-Q2 (previously row 4)
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.03

# Row 6 - 2021-Q1 (previously row 5)
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 4
$total.Range("D6").Value = 1.26

# Row 7 - 2020-Q4 (previously row 6)
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 6
$total.Range("D7").Value = 3.24
